$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing metadata rows ---
$ws.Range("B3").Value = "0.1.7"                         # Version
$ws.Range("B6").Value = "draft"                         # Status
$ws.Range("B8").Value = "2024-08-23T10:17:11-05:00"     # Date

# Row 10: Contact -> new publisher-style contact text
$ws.Range("A10").Value = "Contact"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Row 11: Contact -> Bob Milius
$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Prepare new row 16 with the same formatting as the existing last data row (15)
# before we overwrite/shift row 15's own content, so the border/fill/alignment
# style (s="2") carries through to the newly-added row.
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

# New row 12: Jurisdiction (no value) - everything below shifts down by one row
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# Row 13: Description
$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "CBC W Ordered Manual Differential panel - Blood (57782-5)"

# Row 14: Purpose (no value)
$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").Value = ""

# Row 15: Copyright (no value)
$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").Value = ""

# Row 16: Immutable / BooleanType[null]
$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"
